$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.334129810333252
$ws.Range("B1").Value = 1.418121814727783
$ws.Range("C1").Value = 4.25315523147583
$ws.Range("D1").Value = 3.105960607528687
$ws.Range("E1").Value = 1.013905167579651
